# register.xlsx - "data updated final before launch"
#
# - update the discord-hosted image link text in L2 (shared string) to the
#   freshly uploaded photo's URL
# - turn M2 into a real hyperlink (pointing at the Google-form URL it
#   already displays as text) - Excel auto-creates/assigns the built-in
#   "Hyperlink" cell style (new font + style records) when you do this
# - move the view: scroll column E into view, select L10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- L2: swap the attachment link for the freshly uploaded photo ---------
$ws.Range("L2").Value = "https://media.discordapp.net/attachments/1162451241872412901/1169172195537326130/IMG_20231101_124237.jpg?ex=65546f75&is=6541fa75&hm=ad5ea474b6846acbcb7546d6d3ac8ff47c6d96fc81312f10d6e42e361937ca4b&=&width=733&height=662"

# --- M2: make the existing "Form" URL text a clickable hyperlink ---------
$target = $ws.Range("M2").Value()
[void]$ws.Hyperlinks.Add($ws.Range("M2"), $target)

# --- view state: scroll so column E is left-most, select L10 -------------
[void]$excel.Goto($ws.Range("E1"), $true)
[void]$ws.Range("L10").Select()

Write-Output "done"
